# Attempting to make a dataframe
# Update the sample employee record (row 2) on the AllData sheet and
# move the active selection back to A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tyler"
$ws.Range("C2").Value = "Kathy"

# D2/E2 hold digit-only strings ("42544", "45645343") that must stay text,
# not be auto-converted to numbers - use a leading apostrophe, same as
# typing it directly into Excel.
$ws.Range("D2").Value = "'42544"
$ws.Range("E2").Value = "'45645343"

$ws.Range("A1").Select() | Out-Null
